$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format so numeric-looking / percentage-looking
# strings are preserved exactly as text and not auto-converted by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.336.17'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '1.835.81'
$ws.Range('E3').Value = '  -2.24%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '258.58'
$ws.Range('E5').Value = '  -7.31%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '0.5193'
$ws.Range('E7').Value = '  -1.63%  '
$ws.Range('D8').Value = '0.3224'
$ws.Range('E8').Value = '  -6.46%  '
$ws.Range('D9').Value = '0.06744'
$ws.Range('D10').Value = '18.47'
$ws.Range('E10').Value = '  -8.35%  '
$ws.Range('D11').Value = '0.7580'
$ws.Range('E11').Value = '  -6.16%  '
$ws.Range('D12').Value = '0.07647'
$ws.Range('E12').Value = '  -2.79%  '
$ws.Range('D13').Value = '1.827.49'
$ws.Range('E13').Value = '  -2.79%  '
$ws.Range('D14').Value = '88.25'
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = '5.011'
$ws.Range('E15').Value = '  -3.00%  '
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '13.92'
$ws.Range('E17').Value = '  -4.57%  '
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '0.000007876'
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').Value = '26.358.87'
$ws.Range('E20').Value = '  -2.11%  '
$ws.Range('D21').Value = '2.067.99'
$ws.Range('E21').Value = '  -1.77%  '
$ws.Range('D22').Value = '4.537'
$ws.Range('E22').Value = '  -4.55%  '
$ws.Range('D23').Value = '9.413'
$ws.Range('E23').Value = '  -6.12%  '
$ws.Range('D24').Value = '5.914'
$ws.Range('E24').Value = '  -4.30%  '
$ws.Range('D25').Value = '144.19'
$ws.Range('E25').Value = '  -1.62%  '
$ws.Range('D26').Value = '2.228'
$ws.Range('E26').Value = '  -4.98%  '
$ws.Range('D27').Value = '1.647'
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('D28').Value = '16.87'
$ws.Range('E28').Value = '  -2.90%  '
$ws.Range('D29').Value = '111.26'
$ws.Range('E29').Value = '  -2.36%  '
$ws.Range('E30').Value = '  -5.00%  '
$ws.Range('D31').Value = '4.127'
$ws.Range('E31').Value = '  -4.60%  '
$ws.Range('D32').Value = '0.08698'
$ws.Range('E32').Value = '  -2.52%  '
$ws.Range('D33').Value = '0.04771'
$ws.Range('E33').Value = '  -3.54%  '
$ws.Range('D34').Value = '2.850'
$ws.Range('E35').Value = '  -5.53%  '
$ws.Range('D36').Value = '0.6946'
$ws.Range('E36').Value = '  -5.68%  '
$ws.Range('E37').Value = '  -6.84%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '2.203'
$ws.Range('E38').Value = '  -8.42%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01755'
$ws.Range('E39').Value = '  -5.23%  '
$ws.Range('D40').Value = '0.4825'
$ws.Range('E40').Value = '  -6.25%  '
$ws.Range('D41').Value = '111.16'
$ws.Range('E41').Value = '  -4.40%  '
$ws.Range('D42').Value = '6.102'
$ws.Range('E42').Value = '  -1.75%  '
$ws.Range('D43').Value = '0.8802'
$ws.Range('E43').Value = '  -8.12%  '
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').Value = '7.641'
$ws.Range('E45').Value = '  -5.58%  '
$ws.Range('D46').Value = '0.4124'
$ws.Range('E46').Value = '  -8.66%  '
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').Value = '8.991'
$ws.Range('E48').Value = '  -4.46%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '34.61'
$ws.Range('E49').Value = '  -5.45%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.1218'
$ws.Range('E50').Value = '  -9.42%  '
$ws.Range('D51').Value = '0.8816'
$ws.Range('E51').Value = '  -0.28%  '

# Restore default (Normal) style on the range so no stray number-format
# style is left behind on the cells (matches original formatting).
$ws.Range("D2:E51").Style = "Normal"
